$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 121, shifting all following
# rows (old 121-158) down to 124-161.
$ws.Rows("121:123").Insert()

# Row 121: Murcott / Especial, $/caja 18 kilos, Region Metropolitana
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44553
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100102
$ws.Range("H121").Value = "Cítricos"
$ws.Range("I121").Value = 100102004
$ws.Range("J121").Value = "Mandarina"
$ws.Range("K121").Value = "Murcott"
$ws.Range("L121").Value = "Especial"
$ws.Range("M121").Value = 300
$ws.Range("N121").Value = 16000
$ws.Range("O121").Value = 16000
$ws.Range("P121").Value = 16000
$ws.Range("Q121").Value = "$/caja 18 kilos"
$ws.Range("R121").Value = "Región Metropolitana"
$ws.Range("S121").Value = 889
$ws.Range("T121").Value = 18

# Row 122: Murcott / Primera, $/caja 18 kilos, Region Metropolitana
$ws.Range("A122").Value = 4
$ws.Range("B122").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C122").Value = "Los Lagos"
$ws.Range("D122").Value = 44553
$ws.Range("E122").Value = 10
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100102
$ws.Range("H122").Value = "Cítricos"
$ws.Range("I122").Value = 100102004
$ws.Range("J122").Value = "Mandarina"
$ws.Range("K122").Value = "Murcott"
$ws.Range("L122").Value = "Primera"
$ws.Range("M122").Value = 300
$ws.Range("N122").Value = 14000
$ws.Range("O122").Value = 14000
$ws.Range("P122").Value = 14000
$ws.Range("Q122").Value = "$/caja 18 kilos"
$ws.Range("R122").Value = "Región Metropolitana"
$ws.Range("S122").Value = 778
$ws.Range("T122").Value = 18

# Row 123: Murcott / Segunda, $/caja 18 kilos, Region Metropolitana
$ws.Range("A123").Value = 4
$ws.Range("B123").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C123").Value = "Los Lagos"
$ws.Range("D123").Value = 44553
$ws.Range("E123").Value = 10
$ws.Range("F123").Value = "Fruta"
$ws.Range("G123").Value = 100102
$ws.Range("H123").Value = "Cítricos"
$ws.Range("I123").Value = 100102004
$ws.Range("J123").Value = "Mandarina"
$ws.Range("K123").Value = "Murcott"
$ws.Range("L123").Value = "Segunda"
$ws.Range("M123").Value = 300
$ws.Range("N123").Value = 12000
$ws.Range("O123").Value = 12000
$ws.Range("P123").Value = 12000
$ws.Range("Q123").Value = "$/caja 18 kilos"
$ws.Range("R123").Value = "Región Metropolitana"
$ws.Range("S123").Value = 667
$ws.Range("T123").Value = 18
